$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.608.03"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "3.108.91"
$ws.Range("E3").Value = "  +2.59%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'385.47"
$ws.Range("E5").Value = "  +1.59%  "
$ws.Range("D6").Value = "'103.98"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("E7").Value = "  -1.13%  "
$ws.Range("E9").Value = "  -1.37%  "
$ws.Range("D10").Value = "'37.20"
$ws.Range("E10").Value = "  +1.55%  "
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "3.599.75"
$ws.Range("E13").Value = "  +2.55%  "
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").Value = "'7.85"
$ws.Range("D16").Value = "3.102.00"
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("D18").Value = "'10.92"
$ws.Range("E18").Value = "  +2.93%  "
$ws.Range("D19").Value = "51.626.46"
$ws.Range("D20").Value = "'3.28"
$ws.Range("E20").Value = "  +7.54%  "
$ws.Range("D21").Value = "'12.51"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "'70.03"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "'267.07"
$ws.Range("E24").Value = "  -0.69%  "
$ws.Range("D25").Value = "'3.17"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("E26").Value = "  -1.38%  "
$ws.Range("E27").Value = "  +3.30%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  -6.66%  "
$ws.Range("D30").Value = "'0.167"
$ws.Range("E30").Value = "  -3.41%  "
$ws.Range("E31").Value = "  -1.82%  "
$ws.Range("D32").Value = "'10.43"
$ws.Range("E32").Value = "  +1.59%  "
$ws.Range("D33").Value = "'0.0482"
$ws.Range("E33").Value = "  +6.54%  "
$ws.Range("D34").Value = "'35.28"
$ws.Range("E34").Value = "  +3.11%  "
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("D36").Value = "'50.04"
$ws.Range("E36").Value = "  -1.13%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("D39").Value = "'0.292"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("D41").Value = "'129.22"
$ws.Range("E41").Value = "  +1.71%  "
$ws.Range("D42").Value = "'16.61"
$ws.Range("E42").Value = "  -2.97%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("E44").Value = "  -2.90%  "
$ws.Range("D45").Value = "'3.77"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D46").Value = "'22.28"
$ws.Range("E46").Value = "  +2.12%  "
$ws.Range("E47").Value = "  +6.13%  "
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("D49").Value = "2.071.92"
$ws.Range("E49").Value = "  +1.88%  "
$ws.Range("D50").Value = "'0.949"
$ws.Range("E50").Value = "  +20.68%  "
$ws.Range("D51").Value = "'0.0322"
$ws.Range("E51").Value = "  +0.24%  "
